$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header labels in row 2 (B2 and F2 both become "total")
$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"

# Remove the two section-header-only rows (no data rows):
# Row 8 "grandes regiões e unidades da federação" and Row 5 "situação do domicílio"
# Delete the lower one first so the other row index stays valid.
$ws.Rows(8).Delete()
$ws.Rows(5).Delete()
